$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Rename the existing sheet (USDBTC) to PAIR2 and clear its contents
# so the internal shared-string pool is fully freed, then rebuild it
# (and the new VNDUSD sheet) from scratch in the desired order.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "PAIR2"
$ws1.Cells.Clear()

# Add the second sheet right after PAIR2
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "VNDUSD"
$ws2.Cells.Clear()

# ---- Column widths / formatting (match original layout) ----
foreach ($ws in @($ws1, $ws2)) {
    $ws.Range("E1:F1").ColumnWidth = 35.7109375
    $ws.Range("H1:H1").ColumnWidth = 30.7109375
}

# ---- Sheet1 (PAIR2) header row ----
$ws1.Range("A1").Value = "MONTH"
$ws1.Range("B1").Value = "DAY"
$ws1.Range("C1").Value = "TIME"
$ws1.Range("D1").Value = "POSITION"
$ws1.Range("E1").Value = "15MIN CHART"
$ws1.Range("F1").Value = "1HR CHART"
$ws1.Range("G1").Value = "PROFIT R"
$ws1.Range("H1").Value = "COMMENTS"
$ws1.Range("I1").Value = "ID"
$ws1.Range("J1").Value = "SUM"

# ---- Sheet1 (PAIR2) data row ----
$ws1.Range("A2").Value = 7
$ws1.Range("B2").Value = 6
$ws1.Range("C2").Value = "14:58:01.136710"
$ws1.Range("D2").Value = "Sell"
$ws1.Range("E2").Value = "link 1"
$ws1.Range("F2").Value = "link 2"
$ws1.Range("G2").Value = 3
$ws1.Range("H2").Value = "this is my comment"
$ws1.Range("I2").Value = 806
$ws1.Range("J2").Value = 3

# ---- Sheet2 (VNDUSD) header row ----
$ws2.Range("A1").Value = "MONTH"
$ws2.Range("B1").Value = "DAY"
$ws2.Range("C1").Value = "TIME"
$ws2.Range("D1").Value = "POSITION"
$ws2.Range("E1").Value = "15MIN CHART"
$ws2.Range("F1").Value = "1HR CHART"
$ws2.Range("G1").Value = "PROFIT R"
$ws2.Range("H1").Value = "COMMENTS"
$ws2.Range("I1").Value = "ID"
$ws2.Range("J1").Value = "SUM"

# ---- Sheet2 (VNDUSD) data rows ----
$ws2.Range("A2").Value = 7
$ws2.Range("B2").Value = 6
$ws2.Range("C2").Value = "14:57:30.000793"
$ws2.Range("D2").Value = "Sell"
$ws2.Range("E2").Value = "link 1"
$ws2.Range("F2").Value = "link 2"
$ws2.Range("G2").Value = 3
$ws2.Range("H2").Value = "this is my comment"
$ws2.Range("I2").Value = 606
$ws2.Range("J2").Value = 3

$ws2.Range("A3").Value = 10
$ws2.Range("B3").Value = 6
$ws2.Range("C3").Value = "14:57:03.901108"
$ws2.Range("D3").Value = "Sell"
$ws2.Range("E3").Value = "link 1"
$ws2.Range("F3").Value = "link 2"
$ws2.Range("G3").Value = 3
$ws2.Range("H3").Value = "this is my comment"
$ws2.Range("I3").Value = 406
$ws2.Range("J3").Value = 3

# ---- Apply the border/centered style used by the original sheet (style index 1) ----
foreach ($ws in @($ws1, $ws2)) {
    $used = $ws.Range("A1", $ws.Cells.Item(3, 10))
    $used.Borders.LineStyle = 1
    $used.HorizontalAlignment = -4108
    $used.VerticalAlignment = -4108
}

$ws1.Activate()
